$wb = $excel.ActiveWorkbook

# --- Sheet 1: "basic_information" -> "connected component" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "connected component"

# Set values in the order that matches the target shared-string ordering:
# components first, then the edge-count line.
$ws1.Range("A3").Value = "Component 1: [USAKV, CAYIK, CAYPX]"
$ws1.Range("A4").Value = "Component 2: [USDVO, CAZAM]"
$ws1.Range("A5").Value = "Component 3: [USAMM, CAYUL, USCLT, USDCA, USDFW, USDTW, USJFK, USLAS, USLAX, USLGA, USMCO, USMIA, USORD, USPHL, USATL, USBDL, USBOS, USDEN, USEWR, USFLL, USHOG, USIAD, USIAH, USMSP, USSFO, USALG, CAYYZ, CAYEG, CAYVR, CAYYC, CAYYJ, USSEA, CAYLW, USPDX, USPHX, USSAN, USSLC, USPSP, USHND, USHNL, USSNA, USANU, USBNA, USBWI, USCLE, USCMH, USCVG, USIND, USMCI, USMDT, USMKE, USMSY, USPIT, USRDU, USROC, USRSW, USSTL, USSYR, USTPA, USGEO, USMYR, CAYTZ, USMDW, CAYQR, CAYWG, CAYXE, CAYOW, CAYHZ, CAYQB, CAYQM, CAYYT, CAYMM, CAYKF, CAYXU]"
$ws1.Range("A6").Value = "Component 4: [USLHW, CAYZY]"
$ws1.Range("A2").Value = "total edges is: 170"

$ws1.Columns.Item(1).ColumnWidth = 33.66

$ws1.Range("A11").Select()

# --- Sheet 2: "degree" ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()

$ws2.Range("F:F").Select()

$srt = $ws2.Sort
$srt.SortFields.Clear()
$srt.SortFields.Add($ws2.Range("F1"), 0, 2, 0, 0)
$srt.SetRange($ws2.Range("F1:F83"))
$srt.Header = -4142
$srt.Apply()

Write-Host "done"
